$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GOOGL")

# Row 2 - Revenue
$ws.Range("B2").Value = 196663000000.0
$ws.Range("D2").Value = 171526000000.0
$ws.Range("E2").Value = 165830000000.0
$ws.Range("F2").Value = 166628000000.0
$ws.Range("G2").Value = 161857000000.0

# Row 4 - Gross Profit
$ws.Range("B4").Value = 106810000000.0
$ws.Range("D4").Value = 91854000000.0
$ws.Range("E4").Value = 89707000000.0
$ws.Range("F4").Value = 91762000000.0
$ws.Range("G4").Value = 89961000000.0

# Row 8 - Interest Expense (Operating): was blank inline string, now numeric
$ws.Range("B8").Value = 190000000.0

# Row 9 - Non-operating Income/Expense
$ws.Range("B9").Value = -12419000000.0
$ws.Range("D9").Value = -2465000000.0
$ws.Range("E9").Value = -438000000.0
$ws.Range("F9").Value = 4410000000.0
$ws.Range("G9").Value = 5394000000.0

# Row 10 - Non-operating Interest Expenses: was blank inline string, now numeric
$ws.Range("B10").Value = 1624000000.0

# Row 15 - EPS (Basic)
$ws.Range("B15").Value = 75.12
$ws.Range("D15").Value = 51.89
$ws.Range("E15").Value = 45.69
$ws.Range("F15").Value = 49.89
$ws.Range("G15").Value = 49.6

# Row 16 - Gross Margin
$ws.Range("B16").Value = 0.5431

# Row 17 - EBIT Margin
$ws.Range("B17").Value = 0.2526

# Row 18 - EBT margin
$ws.Range("B18").Value = 0.3133

# Row 19 - Net Profit Margin
$ws.Range("B19").Value = 0.2612

# Row 20 - Free Cash Flow Margin
$ws.Range("B20").Value = 0.1726

# Row 21 - EBITDA
$ws.Range("B21").Value = 63007000000.0
$ws.Range("D21").Value = 48046000000.0
$ws.Range("E21").Value = 45430000000.0
$ws.Range("F21").Value = 47827000000.0
$ws.Range("G21").Value = 46012000000.0

# Row 23 - EPS (Diluted, from Cont. Ops)
$ws.Range("B23").Value = 70.2994

# Row 31 - EBITDA Margin
$ws.Range("B31").Value = 0.3204

# Row 32 - Operating Cash Flow Margin
$ws.Range("B32").Value = 0.371
